$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New file identity: the localization test fixture was regenerated with a new
# guid (a4f23e0e-ac55-41cb-9caa-abf3e936da45) replacing the old one
# (2d391c22-86b7-4001-8730-e8fa8ae5e9ac), plus refreshed handoff content
# hashes/timestamps.
# ---------------------------------------------------------------------------
$oldGuid = "2d391c22-86b7-4001-8730-e8fa8ae5e9ac"
$newGuid = "a4f23e0e-ac55-41cb-9caa-abf3e936da45"
$newZhHash = "66775ff27de1ad9b097a3b1c4f858b3e6700d450"

$newMdName        = "$newGuid.md"
$newMdPath        = "e2e\$newGuid.md"
$newZhXlf         = "$newGuid.$newZhHash.zh-cn.xlf"
$newDeXlf         = "$newGuid.$newZhHash.de-de.xlf"
$newGenerateDate  = "2016-08-31 05:01:34"
$newZhHandoffDate = "2016-08-31 05:01:29"
$clearedDate      = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("G2").Value = $newGenerateDate

$ovB2 = $wsOverview.Range("B2")
$ovAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/527cb0ae0d0681988776dea6535b8e26ce06fe75/e2e/$newGuid.md"
$ovB2.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($ovB2, $ovAddress, "", "", $newMdPath)

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhA2 = $wsZh.Range("A2")
$zhAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/527cb0ae0d0681988776dea6535b8e26ce06fe75/e2e/$newGuid.md"
$zhA2.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($zhA2, $zhAddress, "", "", $newMdName)

$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = $newZhHandoffDate
$wsZh.Range("K2").Value = $clearedDate

# "Latest Target File" (I2) / "Latest Handback File" (J2) are no longer
# populated - handback hasn't produced a target/handback file for this run,
# so the old hyperlink + filenames are removed and the cells go blank.
$zhI2 = $wsZh.Range("I2")
$zhI2.Hyperlinks.Delete()
$zhI2.Style = "Normal"
$zhI2.Value = ""
$wsZh.Range("J2").Value = ""

# ---------------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deA2 = $wsDe.Range("A2")
$deAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/527cb0ae0d0681988776dea6535b8e26ce06fe75/e2e/$newGuid.md"
$deA2.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($deA2, $deAddress, "", "", $newMdName)

$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = $newGenerateDate
$wsDe.Range("K2").Value = $clearedDate

$deI2 = $wsDe.Range("I2")
$deI2.Hyperlinks.Delete()
$deI2.Style = "Normal"
$deI2.Value = ""
$wsDe.Range("J2").Value = ""

# ---------------------------------------------------------------------------
# Column width tweaks on the "zh-cn" / "de-de" sheets: columns I (Latest
# Target File) and J (Latest Handback File) shrink now that they hold short
# / empty values instead of long filenames.
# ---------------------------------------------------------------------------
foreach ($ws in @($wsZh, $wsDe)) {
    $ws.Columns.Item(9).ColumnWidth = 18.6506053379604
    $ws.Columns.Item(10).ColumnWidth = 21.7054770333426
}
